$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unhide columns A, K:L, P:Q (previously hidden helper columns)
$ws.Range("A1").EntireColumn.Hidden = $false
$ws.Range("K1:L1").EntireColumn.Hidden = $false
$ws.Range("P1:Q1").EntireColumn.Hidden = $false

# Rows 57, 58, 59 were using a stray "size 18, no name/color" font (fontId 4).
# Re-apply the formatting already used by equivalent rows elsewhere in the
# sheet (fontId 3: size 18 Arial themed) so the odd font disappears.
$ws.Range("B5:O5").Copy()
$ws.Range("B57:O57").PasteSpecial(-4122)

$ws.Range("B52:O52").Copy()
$ws.Range("B58:O58").PasteSpecial(-4122)

$ws.Range("B3:O3").Copy()
$ws.Range("B59:O59").PasteSpecial(-4122)

# Fill in the remaining-tickets formulas on row 79 (F:H), mirroring the
# "sold" shared formula already present on row 76.
$ws.Range("F79:H79").Formula = "=375-F76"
